# Update "paises" (countries) COVID stats sheet + refresh timestamp.
# The underlying data table (rows 4-220, columns A:H) is sorted by total
# cases (column B) descending. This edit refreshes the case counts for a
# number of countries; several of them (Ecuador/Rumania, Ruanda/Trinidad y
# Tobago, Guyana/Mali, Togo/Republica de Chipre, Nueva Caledonia/Santa
# Lucia, Montserrat/Islas Malvinas) change rank relative to their neighbour
# and so swap rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Octubre de 2020 a las 23:55"

function Set-Row {
    param(
        [int]$Row,
        [string]$Country,
        [double]$Total,
        [double]$Nuevos,
        [double]$Activos,
        [double]$Recuperados,
        [double]$Criticos,
        [double]$MuertesHoy,
        [double]$Muertes
    )
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $Total
    $ws.Cells.Item($Row, 3).Value = $Nuevos
    $ws.Cells.Item($Row, 4).Value = $Activos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $Criticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

Set-Row  4 "Estados Unidos"        7762482 36753 4969405 2576510 0 715 216567
Set-Row  5 "India"                 6832988 78809 5824462  902972 0 963 105554
Set-Row  6 "Brasil"                5000694 29741 4391424  461042 0 657 148228
Set-Row 13 "Sudafrica"              685155  1913  618127   49780 0 145  17248
Set-Row 31 "Ecuador"                143531  1475  120511   11277 0  41  11743
Set-Row 32 "Rumania"                142570  2958  111564   25803 0  82   5203
Set-Row 57 "Barein"                  73932   456   69411    4259 0   0    262
Set-Row 84 "Bulgaria"                22743   437   15448    6422 0  11    873
Set-Row 89 "Costa de Marfil"         19935    32   19550     265 0   0    120
Set-Row 105 "Maldivas"               10656    35    9547    1075 0   0     34
Set-Row 113 "Gabon"                   8815     7    8164     597 0   0     54
Set-Row 121 "Malaui"                  5803     7    4575    1048 0   0    180
Set-Row 131 "Ruanda"                  4883    10    3408    1446 0   0     29
Set-Row 132 "Trinidad yTobago"        4876    30    3010    1782 0   1     84
Set-Row 145 "Guyana"                  3292   104    2084    1113 0   3     95
Set-Row 146 "Mali"                    3210    15    2502     577 0   0    131
Set-Row 153 "Sierra Leona"            2287    10    1716     499 0   0     72
Set-Row 159 "Yemen"                   2049     2    1328     128 0   0    593
Set-Row 160 "Togo"                    1898    17    1419     430 0   0     49
Set-Row 161 "Republica de Chipre"     1897    21    1369     504 0   1     24
Set-Row 176 "Burundi"                  515     1     472      42 0   0      1
Set-Row 207 "Nueva Caledonia"           27     0      27       0 0   0      0
Set-Row 208 "Santa Lucia"               27     0      27       0 0   0      0
Set-Row 215 "Montserrat"                13     0      12       0 0   0      1
Set-Row 216 "Islas Malvinas"            13     0      13       0 0   0      0
